# RHI plots with the position of the sphere added
# Adds two new columns (AA: "Exp Constant", AB: "Exp Constant [dB]") to the
# "tabla" worksheet, mirroring the header style used by the existing table
# header (column Z) and filling the new columns with the constant exponent
# values for every data row (rows 2-34).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("AA1").Value = "Exp Constant"
$ws.Range("AB1").Value = "Exp Constant [dB]"

# Match the formatting (bold, border, centered) used by the rest of the header row
$ws.Range("Z1").Copy()
$ws.Range("AA1:AB1").PasteSpecial(-4122)

# New data columns - same constant value repeated for every row
$ws.Range("AA2:AA34").Value = 385250961.9682089
$ws.Range("AB2:AB34").Value = 85.85743731821252
